$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.1816793333333333
$ws.Cells.Item(2, 8).Value = 0.545038
$ws.Cells.Item(2, 9).Value = 0.3544748187589303
$ws.Cells.Item(2, 10).Value = 0.3544748187589304
$ws.Cells.Item(2, 13).Value = 13.67700833333333
$ws.Cells.Item(2, 14).Value = 41.031025
$ws.Cells.Item(2, 15).Value = 0.124413831206147
$ws.Cells.Item(2, 16).Value = 0.124413831206147
$ws.Cells.Item(2, 17).Value = 2.484829755994444
$ws.Cells.Item(2, 18).Value = 22.36346780395
$ws.Cells.Item(2, 19).Value = 0.04410157026790309
$ws.Cells.Item(2, 20).Value = 0.0441015702679031

$ws.Cells.Item(3, 7).Value = 0.1816793333333333
$ws.Cells.Item(3, 8).Value = 0.545038
$ws.Cells.Item(3, 9).Value = 0.3544748187589303
$ws.Cells.Item(3, 10).Value = 0.3544748187589304
$ws.Cells.Item(3, 13).Value = 74.64939600000001
$ws.Cells.Item(3, 15).Value = 0.679053278848249
$ws.Cells.Item(3, 16).Value = 0.6790532788482488
$ws.Cells.Item(3, 17).Value = 13.562252499016
$ws.Cells.Item(3, 18).Value = 122.060272491144
$ws.Cells.Item(3, 19).Value = 0.2407072879473904
$ws.Cells.Item(3, 20).Value = 0.2407072879473904

$ws.Cells.Item(4, 7).Value = 0.1816793333333333
$ws.Cells.Item(4, 8).Value = 0.545038
$ws.Cells.Item(4, 9).Value = 0.3544748187589303
$ws.Cells.Item(4, 10).Value = 0.3544748187589304
$ws.Cells.Item(4, 13).Value = 1.629335666666667
$ws.Cells.Item(4, 14).Value = 4.888007
$ws.Cells.Item(4, 15).Value = 0.01482136207497777
$ws.Cells.Item(4, 16).Value = 0.01482136207497777
$ws.Cells.Item(4, 17).Value = 0.2960166176962222
$ws.Cells.Item(4, 18).Value = 2.664149559266
$ws.Cells.Item(4, 19).Value = 0.005253799635288229
$ws.Cells.Item(4, 20).Value = 0.00525379963528823

$ws.Cells.Item(5, 7).Value = 0.1816793333333333
$ws.Cells.Item(5, 8).Value = 0.545038
$ws.Cells.Item(5, 9).Value = 0.3544748187589303
$ws.Cells.Item(5, 10).Value = 0.3544748187589304
$ws.Cells.Item(5, 13).Value = 19.17462033333333
$ws.Cells.Item(5, 14).Value = 57.523861
$ws.Cells.Item(5, 15).Value = 0.174423230537864
$ws.Cells.Item(5, 16).Value = 0.174423230537864
$ws.Cells.Item(5, 17).Value = 3.483632239079778
$ws.Cells.Item(5, 18).Value = 31.352690151718
$ws.Cells.Item(5, 19).Value = 0.06182864303225645
$ws.Cells.Item(5, 20).Value = 0.06182864303225647

$ws.Cells.Item(6, 7).Value = 0.1816793333333333
$ws.Cells.Item(6, 8).Value = 0.545038
$ws.Cells.Item(6, 9).Value = 0.3544748187589303
$ws.Cells.Item(6, 10).Value = 0.3544748187589304
$ws.Cells.Item(6, 13).Value = 0.801214
$ws.Cells.Item(6, 14).Value = 2.403642
$ws.Cells.Item(6, 15).Value = 0.007288297332762355
$ws.Cells.Item(6, 16).Value = 0.007288297332762355
$ws.Cells.Item(6, 17).Value = 0.1455640253773333
$ws.Cells.Item(6, 18).Value = 1.310076228396
$ws.Cells.Item(6, 19).Value = 0.002583517876092131
$ws.Cells.Item(6, 20).Value = 0.002583517876092132

$ws.Cells.Item(7, 7).Value = 0.07524966666666667
$ws.Cells.Item(7, 9).Value = 0.14681973708257
$ws.Cells.Item(7, 10).Value = 0.14681973708257
$ws.Cells.Item(7, 13).Value = 13.67700833333333
$ws.Cells.Item(7, 14).Value = 41.031025
$ws.Cells.Item(7, 15).Value = 0.124413831206147
$ws.Cells.Item(7, 16).Value = 0.124413831206147
$ws.Cells.Item(7, 17).Value = 1.029190318080556
$ws.Cells.Item(7, 18).Value = 9.262712862724999
$ws.Cells.Item(7, 19).Value = 0.01826640598712173
$ws.Cells.Item(7, 20).Value = 0.01826640598712173

$ws.Cells.Item(8, 7).Value = 0.07524966666666667
$ws.Cells.Item(8, 9).Value = 0.14681973708257
$ws.Cells.Item(8, 10).Value = 0.14681973708257
$ws.Cells.Item(8, 13).Value = 74.64939600000001
$ws.Cells.Item(8, 15).Value = 0.679053278848249
$ws.Cells.Item(8, 16).Value = 0.6790532788482488
$ws.Cells.Item(8, 17).Value = 5.617342165868001
$ws.Cells.Item(8, 18).Value = 50.556079492812
$ws.Cells.Item(8, 19).Value = 0.09969842386555697
$ws.Cells.Item(8, 20).Value = 0.09969842386555695

$ws.Cells.Item(9, 7).Value = 0.07524966666666667
$ws.Cells.Item(9, 9).Value = 0.14681973708257
$ws.Cells.Item(9, 10).Value = 0.14681973708257
$ws.Cells.Item(9, 13).Value = 1.629335666666667
$ws.Cells.Item(9, 14).Value = 4.888007
$ws.Cells.Item(9, 15).Value = 0.01482136207497777
$ws.Cells.Item(9, 16).Value = 0.01482136207497777
$ws.Cells.Item(9, 17).Value = 0.1226069658047778
$ws.Cells.Item(9, 18).Value = 1.103462692243
$ws.Cells.Item(9, 19).Value = 0.00217606848305381
$ws.Cells.Item(9, 20).Value = 0.00217606848305381

$ws.Cells.Item(10, 7).Value = 0.07524966666666667
$ws.Cells.Item(10, 9).Value = 0.14681973708257
$ws.Cells.Item(10, 10).Value = 0.14681973708257
$ws.Cells.Item(10, 13).Value = 19.17462033333333
$ws.Cells.Item(10, 14).Value = 57.523861
$ws.Cells.Item(10, 15).Value = 0.174423230537864
$ws.Cells.Item(10, 16).Value = 0.174423230537864
$ws.Cells.Item(10, 17).Value = 1.442883788543222
$ws.Cells.Item(10, 18).Value = 12.985954096889
$ws.Cells.Item(10, 19).Value = 0.02560877284866168
$ws.Cells.Item(10, 20).Value = 0.02560877284866168

$ws.Cells.Item(11, 7).Value = 0.07524966666666667
$ws.Cells.Item(11, 9).Value = 0.14681973708257
$ws.Cells.Item(11, 10).Value = 0.14681973708257
$ws.Cells.Item(11, 13).Value = 0.801214
$ws.Cells.Item(11, 14).Value = 2.403642
$ws.Cells.Item(11, 15).Value = 0.007288297332762355
$ws.Cells.Item(11, 16).Value = 0.007288297332762355
$ws.Cells.Item(11, 17).Value = 0.06029108642866667
$ws.Cells.Item(11, 18).Value = 0.542619777858
$ws.Cells.Item(11, 19).Value = 0.001070065898175765
$ws.Cells.Item(11, 20).Value = 0.001070065898175765

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.255602
$ws.Cells.Item(12, 8).Value = 0.766806
$ws.Cells.Item(12, 9).Value = 0.4987054441584996
$ws.Cells.Item(12, 10).Value = 0.4987054441584997
$ws.Cells.Item(12, 13).Value = 13.67700833333333
$ws.Cells.Item(12, 14).Value = 41.031025
$ws.Cells.Item(12, 15).Value = 0.124413831206147
$ws.Cells.Item(12, 16).Value = 0.124413831206147
$ws.Cells.Item(12, 17).Value = 3.495870684016666
$ws.Cells.Item(12, 18).Value = 31.46283615615
$ws.Cells.Item(12, 19).Value = 0.06204585495112212
$ws.Cells.Item(12, 20).Value = 0.06204585495112212

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.255602
$ws.Cells.Item(13, 8).Value = 0.766806
$ws.Cells.Item(13, 9).Value = 0.4987054441584996
$ws.Cells.Item(13, 10).Value = 0.4987054441584997
$ws.Cells.Item(13, 13).Value = 74.64939600000001
$ws.Cells.Item(13, 15).Value = 0.679053278848249
$ws.Cells.Item(13, 16).Value = 0.6790532788482488
$ws.Cells.Item(13, 17).Value = 19.080534916392
$ws.Cells.Item(13, 18).Value = 171.724814247528
$ws.Cells.Item(13, 19).Value = 0.3386475670353015
$ws.Cells.Item(13, 20).Value = 0.3386475670353015

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.255602
$ws.Cells.Item(14, 8).Value = 0.766806
$ws.Cells.Item(14, 9).Value = 0.4987054441584996
$ws.Cells.Item(14, 10).Value = 0.4987054441584997
$ws.Cells.Item(14, 13).Value = 1.629335666666667
$ws.Cells.Item(14, 14).Value = 4.888007
$ws.Cells.Item(14, 15).Value = 0.01482136207497777
$ws.Cells.Item(14, 16).Value = 0.01482136207497777
$ws.Cells.Item(14, 17).Value = 0.4164614550713334
$ws.Cells.Item(14, 18).Value = 3.748153095642
$ws.Cells.Item(14, 19).Value = 0.007391493956635731
$ws.Cells.Item(14, 20).Value = 0.007391493956635731

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.255602
$ws.Cells.Item(15, 8).Value = 0.766806
$ws.Cells.Item(15, 9).Value = 0.4987054441584996
$ws.Cells.Item(15, 10).Value = 0.4987054441584997
$ws.Cells.Item(15, 13).Value = 19.17462033333333
$ws.Cells.Item(15, 14).Value = 57.523861
$ws.Cells.Item(15, 15).Value = 0.174423230537864
$ws.Cells.Item(15, 16).Value = 0.174423230537864
$ws.Cells.Item(15, 17).Value = 4.901071306440667
$ws.Cells.Item(15, 18).Value = 44.109641757966
$ws.Cells.Item(15, 19).Value = 0.08698581465694583
$ws.Cells.Item(15, 20).Value = 0.08698581465694584

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.255602
$ws.Cells.Item(16, 8).Value = 0.766806
$ws.Cells.Item(16, 9).Value = 0.4987054441584996
$ws.Cells.Item(16, 10).Value = 0.4987054441584997
$ws.Cells.Item(16, 13).Value = 0.801214
$ws.Cells.Item(16, 14).Value = 2.403642
$ws.Cells.Item(16, 15).Value = 0.007288297332762355
$ws.Cells.Item(16, 16).Value = 0.007288297332762355
$ws.Cells.Item(16, 17).Value = 0.204791900828
$ws.Cells.Item(16, 18).Value = 1.843127107452
$ws.Cells.Item(16, 19).Value = 0.003634713558494458
$ws.Cells.Item(16, 20).Value = 0.003634713558494459
